$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 49.0
$ws.Range("B1").Value = -5.0
$ws.Range("C1").Value = [double]"1.0E38"
$ws.Range("D1").Value = "HelloOne WorldПривет twoмир!"
